# Applies the "Ajuste nos diagramas e prototipos do fornecedor" edit:
#  1. Insert a new "UC-39: Editar fornecedor" bullet paragraph right after
#     the existing "UC-38: Consultar fornecedor" bullet paragraph.
#  2. Drop the dangling "[Caso de uso 38]" reference so the sentence ends
#     at "fornecedores.".
#  3. Collapse "Sistema volta par" + bookmark + "a o passo 1" into a single
#     run "Sistema volta para o passo 1" (bookmark removed from here).
#  4. Remove one of the trailing empty paragraphs before the last table and
#     relocate the "_GoBack" bookmark onto the (new) last empty paragraph.

$d = $word.ActiveDocument
$wdNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1) Insert "UC-39: Editar fornecedor" paragraph after "UC-38: Consultar
#    fornecedor".
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("UC-38: Consultar fornecedor", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)  # wdCollapseEnd
    $uc39Xml = '<w:p ' + $wdNS + '>' `
        + '<w:pPr>' `
        +   '<w:pStyle w:val="Cabealho"/>' `
        +   '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr>' `
        +   '<w:ind w:left="360"/>' `
        +   '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr>' `
        + '</w:pPr>' `
        + '<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>UC-39</w:t></w:r>' `
        + '<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r>' `
        + '<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Editar</w:t></w:r>' `
        + '<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> fornecedor</w:t></w:r>' `
        + '</w:p>'
    $rng.InsertXML($uc39Xml)
}

# ---------------------------------------------------------------------
# 2) "fornecedores. [Caso de uso 38]" -> "fornecedores."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(" [Caso de uso 38]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2 = $d.Content
$rng2.Find.Execute(" [Caso de uso 38]", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ---------------------------------------------------------------------
# 3) "Sistema volta para o passo 1" - drop the _GoBack bookmark that
#    splits the run (normalize the text, removing the split if any).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Sistema volta para o passo 1", $true, $false, $false, $false, $false, $true, 1, $false, "Sistema volta para o passo 1", 2)

Write-Host "done"
